$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{A="ECs"; B="Cx3cl1"; C="Itgav"; D="ECs"; E=3; F=1; G=7.361448666666667; H=22.084346; I=0.3809728075517136; J=0.3809728075517136; K=3; L=1; M=9.423852333333334; N=28.271557; O=0.06654336290212845; P=0.06654336290212845; Q=69.37320519408023; R=624.3588467467221; S=0.02535121178875642; T=0.02535121178875642},
  @{A="ECs"; B="Cx3cl1"; C="Itgav"; D="FAPs"; E=3; F=1; G=7.361448666666667; H=22.084346; I=0.3809728075517136; J=0.3809728075517136; K=3; L=1; M=50.59256466666667; N=151.777694; O=0.3572423751649123; P=0.3572423751649123; Q=372.4345677086804; R=3351.911109378124; S=0.1360996306430192; T=0.1360996306430192},
  @{A="ECs"; B="Cx3cl1"; C="Itgav"; D="MuSCs"; E=3; F=1; G=7.361448666666667; H=22.084346; I=0.3809728075517136; J=0.3809728075517136; K=3; L=1; M=26.84076266666667; N=80.522288; O=0.1895270158659356; P=0.1895270158659356; Q=197.5868965448498; R=1778.282068903648; S=0.07220463934134365; T=0.07220463934134365},
  @{A="ECs"; B="Cx3cl1"; C="Itgav"; D="Resolving-Mac"; E=3; F=1; G=7.361448666666667; H=22.084346; I=0.3809728075517136; J=0.3809728075517136; K=3; L=1; M=54.762539; N=164.287617; O=0.3866872460670236; P=0.3866872460670236; Q=403.1316197048313; R=3628.184577343482; S=0.1473173257785943; T=0.1473173257785943},
  @{A="FAPs"; B="Cx3cl1"; C="Itgav"; D="ECs"; E=3; F=1; G=10.317205; H=30.951615; I=0.5339403605073807; J=0.5339403605073807; K=3; L=1; M=9.423852333333334; N=28.271557; O=0.06654336290212845; P=0.06654336290212845; Q=97.22781641272834; R=875.050347714555; S=0.03553018717733593; T=0.03553018717733593},
  @{A="FAPs"; B="Cx3cl1"; C="Itgav"; D="FAPs"; E=3; F=1; G=10.317205; H=30.951615; I=0.5339403605073807; J=0.5339403605073807; K=3; L=1; M=50.59256466666667; N=151.777694; O=0.3572423751649123; P=0.3572423751649123; Q=521.9738611417566; R=4697.764750275809; S=0.1907461225840662; T=0.1907461225840662},
  @{A="FAPs"; B="Cx3cl1"; C="Itgav"; D="MuSCs"; E=3; F=1; G=10.317205; H=30.951615; I=0.5339403605073807; J=0.5339403605073807; K=3; L=1; M=26.84076266666667; N=80.522288; O=0.1895270158659356; P=0.1895270158659356; Q=276.9216507883467; R=2492.29485709512; S=0.1011961231773457; T=0.1011961231773457},
  @{A="FAPs"; B="Cx3cl1"; C="Itgav"; D="Resolving-Mac"; E=3; F=1; G=10.317205; H=30.951615; I=0.5339403605073807; J=0.5339403605073807; K=3; L=1; M=54.762539; N=164.287617; O=0.3866872460670236; P=0.3866872460670236; Q=564.9963411834949; R=5084.967070651454; S=0.2064679275686328; T=0.2064679275686328},
  @{A="MuSCs"; B="Cx3cl1"; C="Itgav"; D="ECs"; E=1; F=0.3333333333333333; G=1.634232333333333; H=4.902697; I=0.08457548349701474; J=0.08457548349701474; K=3; L=1; M=9.423852333333334; N=28.271557; O=0.06654336290212845; P=0.06654336290212845; Q=15.40076418769211; R=138.606877689229; S=0.005627937090964828; T=0.005627937090964828},
  @{A="MuSCs"; B="Cx3cl1"; C="Itgav"; D="FAPs"; E=1; F=0.3333333333333333; G=1.634232333333333; H=4.902697; I=0.08457548349701474; J=0.08457548349701474; K=3; L=1; M=50.59256466666667; N=151.777694; O=0.3572423751649123; P=0.3572423751649123; Q=82.68000500452422; R=744.120045040718; S=0.03021394660519439; T=0.03021394660519439},
  @{A="MuSCs"; B="Cx3cl1"; C="Itgav"; D="MuSCs"; E=1; F=0.3333333333333333; G=1.634232333333333; H=4.902697; I=0.08457548349701474; J=0.08457548349701474; K=3; L=1; M=26.84076266666667; N=80.522288; O=0.1895270158659356; P=0.1895270158659356; Q=43.86404220119289; R=394.776379810736; S=0.01602933900260789; T=0.01602933900260789},
  @{A="MuSCs"; B="Cx3cl1"; C="Itgav"; D="Resolving-Mac"; E=1; F=0.3333333333333333; G=1.634232333333333; H=4.902697; I=0.08457548349701474; J=0.08457548349701474; K=3; L=1; M=54.762539; N=164.287617; O=0.3866872460670236; P=0.3866872460670236; Q=89.49471188922766; R=805.4524070030488; S=0.03270426079824764; T=0.03270426079824764},
  @{A="Resolving-Mac"; B="Cx3cl1"; C="Itgav"; D="ECs"; E=1; F=0.3333333333333333; G=0.009880666666666668; H=0.029642; I=0.0005113484438908852; J=0.0005113484438908852; K=3; L=1; M=9.423852333333334; N=28.271557; O=0.06654336290212845; P=0.06654336290212845; Q=0.09311394362155558; R=0.838025492594; S=0.00003402684507126984; T=0.00003402684507126984},
  @{A="Resolving-Mac"; B="Cx3cl1"; C="Itgav"; D="FAPs"; E=1; F=0.3333333333333333; G=0.009880666666666668; H=0.029642; I=0.0005113484438908852; J=0.0005113484438908852; K=3; L=1; M=50.59256466666667; N=151.777694; O=0.3572423751649123; P=0.3572423751649123; Q=0.4998882672831112; R=4.498994405548; S=0.0001826753326324617; T=0.0001826753326324617},
  @{A="Resolving-Mac"; B="Cx3cl1"; C="Itgav"; D="MuSCs"; E=1; F=0.3333333333333333; G=0.009880666666666668; H=0.029642; I=0.0005113484438908852; J=0.0005113484438908852; K=3; L=1; M=26.84076266666667; N=80.522288; O=0.1895270158659356; P=0.1895270158659356; Q=0.2652046289884445; R=2.386841660896; S=0.00009691434463832927; T=0.00009691434463832928},
  @{A="Resolving-Mac"; B="Cx3cl1"; C="Itgav"; D="Resolving-Mac"; E=1; F=0.3333333333333333; G=0.009880666666666668; H=0.029642; I=0.0005113484438908852; J=0.0005113484438908852; K=3; L=1; M=54.762539; N=164.287617; O=0.3866872460670236; P=0.3866872460670236; Q=0.5410903936793333; R=4.869813543114; S=0.0001977319215488244; T=0.0001977319215488244}
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$r = 2
foreach ($row in $rows) {
  $c = 1
  foreach ($colName in $cols) {
    $ws.Cells.Item($r, $c).Value = $row[$colName]
    $c++
  }
  $r++
}

Write-Output ("Rows written: " + ($r - 2))
